$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "tomas2"

$ws.Range("A6").Select()
